# Fruta / hortaliza, semanal
# A new weekly price record was inserted at the top of the Cilantro price
# history (row 76), pushing all subsequent records (old rows 76-164) down
# by one row (new rows 77-165). The new row 76 carries a new date and new
# price figures while reusing the constant attributes (market, region,
# category, unit, origin, etc.) shared by every record in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 76; Excel shifts rows 76-164 down to 77-165
# and copies the row-above formatting (date number format on column D)
# into the newly-inserted, currently blank row.
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new weekly record.
$ws.Cells.Item(76, 1).Value2 = 8
$ws.Cells.Item(76, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(76, 3).Value2 = "Coquimbo"
$ws.Cells.Item(76, 4).Value2 = 44803
$ws.Cells.Item(76, 5).Value2 = 4
$ws.Cells.Item(76, 6).Value2 = 100112040
$ws.Cells.Item(76, 7).Value2 = "Cilantro"
$ws.Cells.Item(76, 8).Value2 = "Sin especificar"
$ws.Cells.Item(76, 9).Value2 = "Primera"
$ws.Cells.Item(76, 10).Value2 = 2000
$ws.Cells.Item(76, 11).Value2 = 2000
$ws.Cells.Item(76, 12).Value2 = 2500
$ws.Cells.Item(76, 13).Value2 = 2250
$ws.Cells.Item(76, 14).Value2 = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(76, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(76, 16).Value2 = 1500
$ws.Cells.Item(76, 17).Value2 = 1.5
$ws.Cells.Item(76, 18).Value2 = "Hortaliza"
